# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled scrape).
# D = Price, E = Volume(1h) change. Both columns are plain text in the sheet
# (e.g. "26.020.41", "0.0" + subscript-5 + "7973"), so D-column updates are
# written with a leading apostrophe to force literal text entry - otherwise COM
# would "smart"-parse numeric-looking strings into numbers/dates and mangle
# things like leading zeros or the subscript digits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '26.020.41'
$ws.Range("E2").Value = '  -0.33%  '

$ws.Range("D3").Value = "'" + '1.644.05'
$ws.Range("E3").Value = '  -1.26%  '

$ws.Range("D4").Value = "'" + '1.004'
$ws.Range("E4").Value = '  +0.15%  '

$ws.Range("D5").Value = "'" + '206.83'
$ws.Range("E5").Value = '  -1.24%  '

$ws.Range("D6").Value = "'" + '0.5162'
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").Value = "'" + '1.004'
$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").Value = "'" + '0.2572'
$ws.Range("E8").Value = '  -2.06%  '

$ws.Range("D9").Value = "'" + '0.06219'
$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("D10").Value = "'" + '20.60'
$ws.Range("E10").Value = '  -2.14%  '

$ws.Range("D11").Value = "'" + '0.07571'
$ws.Range("E11").Value = '  +1.22%  '

$ws.Range("D12").Value = "'" + '1.647.72'
$ws.Range("E12").Value = '  -1.31%  '

$ws.Range("D13").Value = "'" + '4.368'
$ws.Range("E13").Value = '  -0.88%  '

$ws.Range("D14").Value = "'" + '1.872.20'
$ws.Range("E14").Value = '  -0.95%  '

$ws.Range("D15").Value = "'" + '0.5356'
$ws.Range("E15").Value = '  -3.92%  '

$ws.Range("D16").Value = "'" + '0.0₅7973'
$ws.Range("E16").Value = '  +1.50%  '

$ws.Range("D17").Value = "'" + '65.76'
$ws.Range("E17").Value = '  -0.25%  '

$ws.Range("D18").Value = "'" + '26.065.46'
$ws.Range("E18").Value = '  -0.13%  '

$ws.Range("E19").Value = '  +0.06%  '

$ws.Range("D20").Value = "'" + '4.652'
$ws.Range("E20").Value = '  -2.42%  '

$ws.Range("D21").Value = "'" + '186.78'
$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("D22").Value = "'" + '10.02'
$ws.Range("E22").Value = '  -3.12%  '

$ws.Range("D23").Value = "'" + '6.106'
$ws.Range("E23").Value = '  -0.72%  '

$ws.Range("D24").Value = "'" + '1.006'
$ws.Range("E24").Value = '  +0.21%  '

$ws.Range("D25").Value = "'" + '147.54'
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").Value = "'" + '0.1200'
$ws.Range("E26").Value = '  -3.03%  '

$ws.Range("D27").Value = "'" + '7.344'
$ws.Range("E27").Value = '  -2.51%  '

$ws.Range("D28").Value = "'" + '15.55'
$ws.Range("E28").Value = '  -1.78%  '

$ws.Range("D29").Value = "'" + '1.360'
$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").Value = "'" + '0.06020'
$ws.Range("E30").Value = '  -3.46%  '

$ws.Range("D31").Value = "'" + '1.245'
$ws.Range("E31").Value = '  -2.06%  '

$ws.Range("D32").Value = "'" + '3.419'
$ws.Range("E32").Value = '  -1.49%  '

$ws.Range("D33").Value = "'" + '3.381'
$ws.Range("E33").Value = '  -0.82%  '

$ws.Range("D34").Value = "'" + '1.616'
$ws.Range("E34").Value = '  -0.22%  '

$ws.Range("D35").Value = "'" + '0.9681'
$ws.Range("E35").Value = '  -2.58%  '

$ws.Range("D36").Value = "'" + '2.384'
$ws.Range("E36").Value = '  -0.89%  '

$ws.Range("D37").Value = "'" + '2.729'
$ws.Range("E37").Value = '  +1.15%  '

$ws.Range("D38").Value = "'" + '0.5860'
$ws.Range("E38").Value = '  -2.48%  '

$ws.Range("D39").Value = "'" + '0.01593'
$ws.Range("E39").Value = '  -0.52%  '

$ws.Range("D40").Value = "'" + '1.072.19'
$ws.Range("E40").Value = '  -0.05%  '

$ws.Range("D41").Value = "'" + '5.842'
$ws.Range("E41").Value = '  -4.54%  '

$ws.Range("D42").Value = "'" + '0.8452'
$ws.Range("E42").Value = '  -1.71%  '

$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("D44").Value = "'" + '100.41'
$ws.Range("E44").Value = '  +1.47%  '

$ws.Range("D45").Value = "'" + '1.802.27'
$ws.Range("E45").Value = '  -0.59%  '

$ws.Range("D46").Value = "'" + '0.0₈108'
$ws.Range("E46").Value = '  -2.37%  '

$ws.Range("D47").Value = "'" + '0.9972'
$ws.Range("E47").Value = '  -0.60%  '

$ws.Range("D48").Value = "'" + '54.39'
$ws.Range("E48").Value = '  -2.61%  '

$ws.Range("D49").Value = "'" + '7.986'
$ws.Range("E49").Value = '  +0.84%  '

$ws.Range("D50").Value = "'" + '0.05210'
$ws.Range("E50").Value = '  -0.73%  '

$ws.Range("D51").Value = "'" + '0.4242'
$ws.Range("E51").Value = '  -0.19%  '
